$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

# Add the new DEFAULT_SEARCH_FOLDER configuration row
$ws.Range("A7").Value = "DEFAULT_SEARCH_FOLDER"
$ws.Range("B7").Value = "C:/Users/USER/National University of Singapore/MTech Y1S2 - Documents/General/Final Capstone Project/Dataset"

# Give the new key cell (A7) a left/right outline so it lines up visually with the boxed table above it
$ws.Range("A7").Borders.Item(7).LineStyle = 1
$ws.Range("A7").Borders.Item(7).Weight = 2
$ws.Range("A7").Borders.Item(10).LineStyle = 1
$ws.Range("A7").Borders.Item(10).Weight = 2

# Draw a full thin box border around the header row and the existing key/value rows
$ws.Range("A1:B6").Borders.LineStyle = 1
$ws.Range("A1:B6").Borders.Weight = 2

# Widen column B so the longer folder path value fits
$ws.Columns.Item(2).ColumnWidth = 39.5703125

# Match the updated selection / active cell from the authored workbook
$ws.Range("D12").Select()

# Explicitly set the page to portrait orientation
$ws.PageSetup.Orientation = 1
